$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 3745
$endRow = 3840
$numRows = $endRow - $startRow + 1

$data = New-Object 'object[,]' $numRows,3

$data[0,0] = 46060.95833333334
$data[0,1] = 46061
$data[0,2] = 108.87
$data[1,0] = 46060.96875
$data[1,1] = 46061.01041666666
$data[1,2] = 101.13
$data[2,0] = 46060.97916666666
$data[2,1] = 46061.02083333334
$data[2,2] = 97.81999999999999
$data[3,0] = 46060.98958333334
$data[3,1] = 46061.03125
$data[3,2] = 93.06999999999999
$data[4,0] = 46061
$data[4,1] = 46061.04166666666
$data[4,2] = 97.79000000000001
$data[5,0] = 46061.01041666666
$data[5,1] = 46061.05208333334
$data[5,2] = 98.66
$data[6,0] = 46061.02083333334
$data[6,1] = 46061.0625
$data[6,2] = 97.09999999999999
$data[7,0] = 46061.03125
$data[7,1] = 46061.07291666666
$data[7,2] = 95.20999999999999
$data[8,0] = 46061.04166666666
$data[8,1] = 46061.08333333334
$data[8,2] = 98.05
$data[9,0] = 46061.05208333334
$data[9,1] = 46061.09375
$data[9,2] = 95.23
$data[10,0] = 46061.0625
$data[10,1] = 46061.10416666666
$data[10,2] = 92.69
$data[11,0] = 46061.07291666666
$data[11,1] = 46061.11458333334
$data[11,2] = 91.63
$data[12,0] = 46061.08333333334
$data[12,1] = 46061.125
$data[12,2] = 93.48999999999999
$data[13,0] = 46061.09375
$data[13,1] = 46061.13541666666
$data[13,2] = 92.25
$data[14,0] = 46061.10416666666
$data[14,1] = 46061.14583333334
$data[14,2] = 91.13
$data[15,0] = 46061.11458333334
$data[15,1] = 46061.15625
$data[15,2] = 91.61
$data[16,0] = 46061.125
$data[16,1] = 46061.16666666666
$data[16,2] = 90.34999999999999
$data[17,0] = 46061.13541666666
$data[17,1] = 46061.17708333334
$data[17,2] = 90.43000000000001
$data[18,0] = 46061.14583333334
$data[18,1] = 46061.1875
$data[18,2] = 90.73
$data[19,0] = 46061.15625
$data[19,1] = 46061.19791666666
$data[19,2] = 90.33
$data[20,0] = 46061.16666666666
$data[20,1] = 46061.20833333334
$data[20,2] = 92.12
$data[21,0] = 46061.17708333334
$data[21,1] = 46061.21875
$data[21,2] = 91.44
$data[22,0] = 46061.1875
$data[22,1] = 46061.22916666666
$data[22,2] = 90.36
$data[23,0] = 46061.19791666666
$data[23,1] = 46061.23958333334
$data[23,2] = 91.61
$data[24,0] = 46061.20833333334
$data[24,1] = 46061.25
$data[24,2] = 84.8
$data[25,0] = 46061.21875
$data[25,1] = 46061.26041666666
$data[25,2] = 87.67
$data[26,0] = 46061.22916666666
$data[26,1] = 46061.27083333334
$data[26,2] = 90.98
$data[27,0] = 46061.23958333334
$data[27,1] = 46061.28125
$data[27,2] = 95.95999999999999
$data[28,0] = 46061.25
$data[28,1] = 46061.29166666666
$data[28,2] = 88.2
$data[29,0] = 46061.26041666666
$data[29,1] = 46061.30208333334
$data[29,2] = 91.33
$data[30,0] = 46061.27083333334
$data[30,1] = 46061.3125
$data[30,2] = 94.69
$data[31,0] = 46061.28125
$data[31,1] = 46061.32291666666
$data[31,2] = 99.09999999999999
$data[32,0] = 46061.29166666666
$data[32,1] = 46061.33333333334
$data[32,2] = 94.41
$data[33,0] = 46061.30208333334
$data[33,1] = 46061.34375
$data[33,2] = 98.72
$data[34,0] = 46061.3125
$data[34,1] = 46061.35416666666
$data[34,2] = 100.1
$data[35,0] = 46061.32291666666
$data[35,1] = 46061.36458333334
$data[35,2] = 100.04
$data[36,0] = 46061.33333333334
$data[36,1] = 46061.375
$data[36,2] = 100.16
$data[37,0] = 46061.34375
$data[37,1] = 46061.38541666666
$data[37,2] = 99.87
$data[38,0] = 46061.35416666666
$data[38,1] = 46061.39583333334
$data[38,2] = 99.48999999999999
$data[39,0] = 46061.36458333334
$data[39,1] = 46061.40625
$data[39,2] = 93.78
$data[40,0] = 46061.375
$data[40,1] = 46061.41666666666
$data[40,2] = 100.9
$data[41,0] = 46061.38541666666
$data[41,1] = 46061.42708333334
$data[41,2] = 98
$data[42,0] = 46061.39583333334
$data[42,1] = 46061.4375
$data[42,2] = 94.86
$data[43,0] = 46061.40625
$data[43,1] = 46061.44791666666
$data[43,2] = 93.12
$data[44,0] = 46061.41666666666
$data[44,1] = 46061.45833333334
$data[44,2] = 96.90000000000001
$data[45,0] = 46061.42708333334
$data[45,1] = 46061.46875
$data[45,2] = 93.98
$data[46,0] = 46061.4375
$data[46,1] = 46061.47916666666
$data[46,2] = 92.53
$data[47,0] = 46061.44791666666
$data[47,1] = 46061.48958333334
$data[47,2] = 90.44
$data[48,0] = 46061.45833333334
$data[48,1] = 46061.5
$data[48,2] = 92.03
$data[49,0] = 46061.46875
$data[49,1] = 46061.51041666666
$data[49,2] = 90.5
$data[50,0] = 46061.47916666666
$data[50,1] = 46061.52083333334
$data[50,2] = 92.56
$data[51,0] = 46061.48958333334
$data[51,1] = 46061.53125
$data[51,2] = 89.59
$data[52,0] = 46061.5
$data[52,1] = 46061.54166666666
$data[52,2] = 90.8
$data[53,0] = 46061.51041666666
$data[53,1] = 46061.55208333334
$data[53,2] = 88.97
$data[54,0] = 46061.52083333334
$data[54,1] = 46061.5625
$data[54,2] = 90.88
$data[55,0] = 46061.53125
$data[55,1] = 46061.57291666666
$data[55,2] = 89.17
$data[56,0] = 46061.54166666666
$data[56,1] = 46061.58333333334
$data[56,2] = 88.76000000000001
$data[57,0] = 46061.55208333334
$data[57,1] = 46061.59375
$data[57,2] = 89.90000000000001
$data[58,0] = 46061.5625
$data[58,1] = 46061.60416666666
$data[58,2] = 91.29000000000001
$data[59,0] = 46061.57291666666
$data[59,1] = 46061.61458333334
$data[59,2] = 93.73
$data[60,0] = 46061.58333333334
$data[60,1] = 46061.625
$data[60,2] = 92.69
$data[61,0] = 46061.59375
$data[61,1] = 46061.63541666666
$data[61,2] = 96.51000000000001
$data[62,0] = 46061.60416666666
$data[62,1] = 46061.64583333334
$data[62,2] = 99.09999999999999
$data[63,0] = 46061.61458333334
$data[63,1] = 46061.65625
$data[63,2] = 105.63
$data[64,0] = 46061.625
$data[64,1] = 46061.66666666666
$data[64,2] = 92.94
$data[65,0] = 46061.63541666666
$data[65,1] = 46061.67708333334
$data[65,2] = 100.57
$data[66,0] = 46061.64583333334
$data[66,1] = 46061.6875
$data[66,2] = 107.46
$data[67,0] = 46061.65625
$data[67,1] = 46061.69791666666
$data[67,2] = 118.19
$data[68,0] = 46061.66666666666
$data[68,1] = 46061.70833333334
$data[68,2] = 111.91
$data[69,0] = 46061.67708333334
$data[69,1] = 46061.71875
$data[69,2] = 120.09
$data[70,0] = 46061.6875
$data[70,1] = 46061.72916666666
$data[70,2] = 125.41
$data[71,0] = 46061.69791666666
$data[71,1] = 46061.73958333334
$data[71,2] = 129.3
$data[72,0] = 46061.70833333334
$data[72,1] = 46061.75
$data[72,2] = 126.95
$data[73,0] = 46061.71875
$data[73,1] = 46061.76041666666
$data[73,2] = 129.89
$data[74,0] = 46061.72916666666
$data[74,1] = 46061.77083333334
$data[74,2] = 132.55
$data[75,0] = 46061.73958333334
$data[75,1] = 46061.78125
$data[75,2] = 129.84
$data[76,0] = 46061.75
$data[76,1] = 46061.79166666666
$data[76,2] = 132.96
$data[77,0] = 46061.76041666666
$data[77,1] = 46061.80208333334
$data[77,2] = 129.55
$data[78,0] = 46061.77083333334
$data[78,1] = 46061.8125
$data[78,2] = 125.09
$data[79,0] = 46061.78125
$data[79,1] = 46061.82291666666
$data[79,2] = 120
$data[80,0] = 46061.79166666666
$data[80,1] = 46061.83333333334
$data[80,2] = 122.21
$data[81,0] = 46061.80208333334
$data[81,1] = 46061.84375
$data[81,2] = 113.65
$data[82,0] = 46061.8125
$data[82,1] = 46061.85416666666
$data[82,2] = 109.35
$data[83,0] = 46061.82291666666
$data[83,1] = 46061.86458333334
$data[83,2] = 103.88
$data[84,0] = 46061.83333333334
$data[84,1] = 46061.875
$data[84,2] = 109.68
$data[85,0] = 46061.84375
$data[85,1] = 46061.88541666666
$data[85,2] = 107.09
$data[86,0] = 46061.85416666666
$data[86,1] = 46061.89583333334
$data[86,2] = 103.63
$data[87,0] = 46061.86458333334
$data[87,1] = 46061.90625
$data[87,2] = 101.05
$data[88,0] = 46061.875
$data[88,1] = 46061.91666666666
$data[88,2] = 105.51
$data[89,0] = 46061.88541666666
$data[89,1] = 46061.92708333334
$data[89,2] = 104.16
$data[90,0] = 46061.89583333334
$data[90,1] = 46061.9375
$data[90,2] = 104.14
$data[91,0] = 46061.90625
$data[91,1] = 46061.94791666666
$data[91,2] = 101.42
$data[92,0] = 46061.91666666666
$data[92,1] = 46061.95833333334
$data[92,2] = 102.25
$data[93,0] = 46061.92708333334
$data[93,1] = 46061.96875
$data[93,2] = 101.71
$data[94,0] = 46061.9375
$data[94,1] = 46061.97916666666
$data[94,2] = 100.28
$data[95,0] = 46061.94791666666
$data[95,1] = 46061.98958333334
$data[95,2] = 93.73999999999999

$targetRange = $ws.Range("A" + $startRow + ":C" + $endRow)
$targetRange.Value = $data

$dateRange = $ws.Range("A" + $startRow + ":B" + $endRow)
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"
